$d = $word.ActiveDocument

# Locate the exact phrase "[P.TIPO_ABR] no.: [P.NUM]/[P.AN]" inside the document
# (the header line that identifies the proceeding type/number/year) and
# replace it with "[P.TIPO_ABR] nº: [P.NUM]/[P.ANO]", split across runs the
# same way the reference edit was produced (first "no." -> "nº", then
# "AN]" -> "ANO]").
$findRng = $d.Content
$found = $findRng.Find.Execute("[P.TIPO_ABR] no.: [P.NUM]/[P.AN]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target placeholder text '[P.TIPO_ABR] no.: [P.NUM]/[P.AN]'"
}

# Re-materialize a plain Range over the same span; re-using the Find-chained
# range object directly causes InsertXML to append instead of replace.
$rng = $d.Range($findRng.Start, $findRng.End)

$rPr = '<w:rPr><w:rFonts w:cs="Cambria" w:ascii="Cambria" w:hAnsi="Cambria"/><w:b/></w:rPr>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r>' + $rPr + '<w:t>[P.TIPO_ABR] n</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>' + [char]0x00BA + '</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>: [P.NUM]/[P.AN</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>O</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>]</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$rng.InsertXML($xml)
